$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Text constants involved in the edit.
# ---------------------------------------------------------------------
$oldCombinedText = "In regard to reproducibility and replicability of our work, a repository is already available containing all the methods and data used in the current work. This repository is currently private in wait for the status of the current submission or the request by the Journal. In addition to this, we are now aiming to prepare a submission for Data in Brief, as we highly encourage and support the reproducibility of academic work. "

$newIntroText = "In regard to reproducibility and replicability of our work, a repository is already available containing all the methods and data used in the current work. This repository is currently private in wait for the status of the current submission or the request by the Journal. "

$newSecondText = "Our current work is part of a general methodology that we are studying for calibration maintenance, and so we aim to prepare a MethodsX paper in the near future gathering the methodology for this purpose. As mentioned in our first submission, the option for Data in Brief is currently under consideration of the researchers who directly made the measurements."

# ---------------------------------------------------------------------
# Locate, by scanning the paragraph collection (not by fixed index, so
# the script is resilient to any paragraph count differences):
#   - the empty paragraph right before the one holding the old
#     combined "reproducibility / Data in Brief" text;
#   - the paragraph holding that old combined text;
#   - the "Closing" style paragraph made of two runs "V" + "aleria ...".
# ---------------------------------------------------------------------
$emptyParaIdx = -1
$oldTextParaIdx = -1
$nameParaIdx = -1

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -eq ($oldCombinedText + "`r")) {
        $oldTextParaIdx = $i
        $emptyParaIdx = $i - 1
    }
    if ($t -eq "Valeria Fonseca Diaz`r" -and $d.Paragraphs.Item($i).Style.NameLocal -eq "Closing") {
        $nameParaIdx = $i
    }
}

# ---------------------------------------------------------------------
# 1) Fill the empty paragraph with the intro text (first half of the
#    text that used to be combined in the following paragraph).
#    A direct Range.Text assignment keeps the paragraph's original
#    (empty) run-properties element intact.
# ---------------------------------------------------------------------
$emptyPara = $d.Paragraphs.Item($emptyParaIdx)
$emptyPara.Range.Text = $newIntroText

# ---------------------------------------------------------------------
# 2) Replace the old combined text with the new MethodsX paragraph.
#    Re-using the same run via Range.Text= keeps xml:space="preserve"
#    stuck from the old text, so instead we splice in a brand-new
#    paragraph (inheriting style/formatting) immediately before the
#    old one, populate it, then delete the old paragraph outright.
# ---------------------------------------------------------------------
$oldTextPara = $d.Paragraphs.Item($oldTextParaIdx)
$r = $oldTextPara.Range
$r.Collapse(1)              # wdCollapseStart
$r.InsertParagraphBefore()
$newTextPara = $d.Paragraphs.Item($oldTextParaIdx)
$newTextPara.Range.Text = $newSecondText
$staleTextPara = $d.Paragraphs.Item($oldTextParaIdx + 1)
$staleTextPara.Range.Delete()

# ---------------------------------------------------------------------
# 3) Merge the "V" + "aleria Fonseca Diaz" runs (Closing-style
#    signature paragraph) into a single run "Valeria Fonseca Diaz".
#    Same splice technique, since the paragraph holds two runs and a
#    plain Range.Text= only overwrites the first of the two.
# ---------------------------------------------------------------------
# Re-scan in case indices shifted due to the paragraph splice above.
$nameParaIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -eq "Valeria Fonseca Diaz`r" -and $d.Paragraphs.Item($i).Style.NameLocal -eq "Closing") {
        $nameParaIdx = $i
    }
}

$namePara = $d.Paragraphs.Item($nameParaIdx)
$rn = $namePara.Range
$rn.Collapse(1)
$rn.InsertParagraphBefore()
$newNamePara = $d.Paragraphs.Item($nameParaIdx)
$newNamePara.Range.Text = "Valeria Fonseca Diaz"
$staleNamePara = $d.Paragraphs.Item($nameParaIdx + 1)
$staleNamePara.Range.Delete()
